$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reaction15")

# Clear out the old range of values (A1:R1) first, then set the new ones.
$ws.Range("A1:R1").ClearContents()

$ws.Range("A1").Value = 30
$ws.Range("B1").Value = 31
